# Fix demo cases and warnings
# - Remove now-unused columns G:H (monetary-related extra output columns) from the
#   key_outputs sheet for rows 2-8.
# - Update the sheet's active cell/selection from F6 to F2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("key_outputs")
$ws.Activate()

# Delete the contents of columns G and H for data rows 2-8 (clears the cells,
# matching the removal of these <c> elements from the sheet XML).
$ws.Range("G2:H8").ClearContents()

# Update the active selection to F2 (previously F6).
$ws.Range("F2").Select()

$wb.Save()
